$d = $word.ActiveDocument

# 1) "recebera" -> "receberá"
$d.Content.Find.Execute("recebera 10 cartas", $true, $false, $false, $false, $false, $true, 1, $false, "receberá 10 cartas", 2)

# 2) "ira receber" -> "irá receber"
$d.Content.Find.Execute("ira receber", $true, $false, $false, $false, $false, $true, 1, $false, "irá receber", 2)

# 3) "dentro de estas 10 cartas existe a possibilidade de entre essas 10 cartas existirem cartas especiais estas cartas especiais" -> "dentro destas 10 cartas existe a possibilidade de existirem cartas especiais, estas cartas especiais"
$d.Content.Find.Execute("dentro de estas 10 cartas existe a possibilidade de entre essas 10 cartas existirem cartas especiais estas cartas especiais", $true, $false, $false, $false, $false, $true, 1, $false, "dentro destas 10 cartas existe a possibilidade de existirem cartas especiais, estas cartas especiais", 2)

# 4) "adversaria." -> "adversária."
$d.Content.Find.Execute("adversaria.", $true, $false, $false, $false, $false, $true, 1, $false, "adversária.", 2)

# 5) "(vermelho, amarelo e laranja) estas cores" -> "(vermelho, amarelo e laranja), estas cores"
$d.Content.Find.Execute("(vermelho, amarelo e laranja) estas cores", $true, $false, $false, $false, $false, $true, 1, $false, "(vermelho, amarelo e laranja), estas cores", 2)

# 6) "influencia o jogo" -> "influência o jogo"
$d.Content.Find.Execute("influencia o jogo", $true, $false, $false, $false, $false, $true, 1, $false, "influência o jogo", 2)
